# Add data for 2021-11-18: refresh the "through" date from Nov 09 to Nov 10
# and bump the affected neighborhood/month carjacking counts accordingly
# (the current running month plus a handful of back-filled prior-year cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Sheet title / rolling header label -------------------------------------------------
$ws.Name = "Through 2021-11-10"
$ws.Range("B1").Value = "November 2021 (through November 10)"

# --- Updated counts (existing values incremented) ---------------------------------------
$ws.Range("B3").Value   = 2   # Garfield Park, Nov 2021 (through Nov 10)
$ws.Range("X3").Value   = 3   # Garfield Park, Nov 2019
$ws.Range("M4").Value   = 4   # Austin, Nov 2020
$ws.Range("AT4").Value  = 6   # Austin, Nov 2017
$ws.Range("AT6").Value  = 2   # West Town, Nov 2017
$ws.Range("M8").Value   = 3   # South Shore, Nov 2020
$ws.Range("B15").Value  = 2   # Grand Crossing, Nov 2021 (through Nov 10)
$ws.Range("M16").Value  = 3   # Washington Heights, Nov 2020
$ws.Range("X16").Value  = 2   # Washington Heights, Nov 2019
$ws.Range("AI23").Value = 2   # Avondale, Nov 2018
$ws.Range("M31").Value  = 2   # Albany Park, Nov 2020
$ws.Range("AT32").Value = 2   # Little Italy, UIC, Nov 2017
$ws.Range("B47").Value  = 3   # Roseland, Nov 2021 (through Nov 10)

# --- Newly populated cells (previously blank) --------------------------------------------
$ws.Range("AT3").Value  = 1   # Garfield Park, Nov 2017
$ws.Range("BE6").Value  = 1   # West Town, Nov 2016
$ws.Range("AI16").Value = 1   # Washington Heights, Nov 2018
$ws.Range("M18").Value  = 1   # Calumet Heights, Nov 2020
$ws.Range("X72").Value  = 1   # Gage Park, Nov 2019
$ws.Range("X84").Value  = 1   # Morgan Park, Nov 2019
$ws.Range("B98").Value  = 1   # Woodlawn, Nov 2021 (through Nov 10)
$ws.Range("AT98").Value = 1   # Woodlawn, Nov 2017
